# Insert a new weekly data row for "Macroferia Regional de Talca - Apio"
# at row 379, pushing the existing rows 379:402 down to 380:403 and
# expanding the sheet dimension to A1:R403.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 379 (shifts rows 379-402 down to 380-403,
# carrying cell formatting such as the date style on column D along with them).
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row 379 with the new record.
$ws.Range("A379").Value = 5
$ws.Range("B379").Value = "Macroferia Regional de Talca"
$ws.Range("C379").Value = "Maule"
$ws.Range("D379").Value = 45265
$ws.Range("E379").Value = 7
$ws.Range("F379").Value = 100112017
$ws.Range("G379").Value = "Apio"
$ws.Range("H379").Value = "Americana (o)"
$ws.Range("I379").Value = "Primera"
$ws.Range("J379").Value = 300
$ws.Range("K379").Value = 15000
$ws.Range("L379").Value = 15000
$ws.Range("M379").Value = 15000
$ws.Range("N379").Value = "`$/docena de matas"
$ws.Range("O379").Value = "Provincia del Elquí"
$ws.Range("P379").Value = 2500
$ws.Range("Q379").Value = 6
$ws.Range("R379").Value = "Hortaliza"
